$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.891.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.37%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4650"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3712"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07379"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8750"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.46"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.773.74"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.372"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.503"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07049"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.911.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.077.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.910"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.48"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.144"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.289"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08930"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7575"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.159"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.461"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.912"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.0000"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.107"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.32%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.30%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.403"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.58%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.924"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.242"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5296"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.52%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5001"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.33"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.666"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06295"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.61%  "
